$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.Value = "'" + $Text
}

Set-TextValue $ws.Range("D2") "283.44"
Set-TextValue $ws.Range("E2") "1.96%"
Set-TextValue $ws.Range("G2") "18"

Set-TextValue $ws.Range("D3") "28.29"
Set-TextValue $ws.Range("E3") "3.53%"
Set-TextValue $ws.Range("G3") "18"

Set-TextValue $ws.Range("D4") "5.079"
Set-TextValue $ws.Range("E4") "4.96%"
Set-TextValue $ws.Range("G4") "18"

Set-TextValue $ws.Range("D5") "0.06642"
Set-TextValue $ws.Range("E5") "4.08%"
Set-TextValue $ws.Range("G5") "18"

Set-TextValue $ws.Range("D6") "7.296"
Set-TextValue $ws.Range("E6") "3.95%"
Set-TextValue $ws.Range("G6") "18"

Set-TextValue $ws.Range("B7") "GateToken"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D7") "3.370"
Set-TextValue $ws.Range("E7") "2.28%"
Set-TextValue $ws.Range("G7") "18"

Set-TextValue $ws.Range("B8") "FTXToken"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D8") "1.365"
Set-TextValue $ws.Range("E8") "3.69%"
Set-TextValue $ws.Range("G8") "18"

Set-TextValue $ws.Range("B9") "MXToken"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D9") "0.9329"
Set-TextValue $ws.Range("E9") "4.83%"
Set-TextValue $ws.Range("G9") "18"

Set-TextValue $ws.Range("B10") "WazirX"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1568"
Set-TextValue $ws.Range("E10") "3.26%"
Set-TextValue $ws.Range("G10") "18"

Set-TextValue $ws.Range("B11") "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D11") "0.06241"
Set-TextValue $ws.Range("E11") "14.24%"
Set-TextValue $ws.Range("G11") "18"

Set-TextValue $ws.Range("B12") "MandalaExchangeToken"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.07526"
Set-TextValue $ws.Range("E12") "1.11%"
Set-TextValue $ws.Range("G12") "18"

Set-TextValue $ws.Range("B13") "BitrueCoin"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.02941"
Set-TextValue $ws.Range("E13") "1.69%"
Set-TextValue $ws.Range("G13") "18"

Set-TextValue $ws.Range("B14") "BitMartToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.08929"
Set-TextValue $ws.Range("E14") "-0.28%"
Set-TextValue $ws.Range("G14") "18"

Set-TextValue $ws.Range("B15") "BitForexToken"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001560"
Set-TextValue $ws.Range("E15") "-1.30%"
Set-TextValue $ws.Range("G15") "18"

Set-TextValue $ws.Range("B16") "CoinExToken"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D16") "0.04449"
Set-TextValue $ws.Range("E16") "1.33%"
Set-TextValue $ws.Range("G16") "18"

Set-TextValue $ws.Range("B17") "One"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D17") "0.0006398"
Set-TextValue $ws.Range("E17") "0.63%"
Set-TextValue $ws.Range("G17") "18"

Set-TextValue $ws.Range("B18") "TigerCash"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D18") "0.006129"
Set-TextValue $ws.Range("E18") "1.86%"
Set-TextValue $ws.Range("G18") "18"

Set-TextValue $ws.Range("B19") "LEO"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D19") "3.468"
Set-TextValue $ws.Range("E19") "-0.13%"
Set-TextValue $ws.Range("G19") "18"

Set-TextValue $ws.Range("D20") "2.233"
Set-TextValue $ws.Range("E20") "-0.05%"
Set-TextValue $ws.Range("G20") "18"

Set-TextValue $ws.Range("D21") "0.3195"
Set-TextValue $ws.Range("E21") "0.79%"
Set-TextValue $ws.Range("G21") "18"

Set-TextValue $ws.Range("D22") "0.1298"
Set-TextValue $ws.Range("E22") "-3.78%"
Set-TextValue $ws.Range("G22") "18"

Set-TextValue $ws.Range("D23") "4.045"
Set-TextValue $ws.Range("E23") "3.14%"
Set-TextValue $ws.Range("G23") "18"

Set-TextValue $ws.Range("D24") "0.1514"
Set-TextValue $ws.Range("E24") "0.57%"
Set-TextValue $ws.Range("G24") "18"

Set-TextValue $ws.Range("D25") "0.001175"
Set-TextValue $ws.Range("E25") "0.18%"
Set-TextValue $ws.Range("G25") "18"

Set-TextValue $ws.Range("D26") "0.004453"
Set-TextValue $ws.Range("E26") "3.95%"
Set-TextValue $ws.Range("G26") "18"

Set-TextValue $ws.Range("D27") "0.0001240"
Set-TextValue $ws.Range("E27") "5.35%"
Set-TextValue $ws.Range("G27") "18"

Set-TextValue $ws.Range("D28") "0.0001605"
Set-TextValue $ws.Range("E28") "-16.88%"
Set-TextValue $ws.Range("G28") "18"

Set-TextValue $ws.Range("G29") "18"

Set-TextValue $ws.Range("G30") "18"

Set-TextValue $ws.Range("G31") "18"

Set-TextValue $ws.Range("G32") "18"

Set-TextValue $ws.Range("G33") "18"

Set-TextValue $ws.Range("G34") "18"

Set-TextValue $ws.Range("G35") "18"

Set-TextValue $ws.Range("G36") "18"

Set-TextValue $ws.Range("G37") "18"

Set-TextValue $ws.Range("G38") "18"

Set-TextValue $ws.Range("G39") "18"

Set-TextValue $ws.Range("D40") "0.04159"
Set-TextValue $ws.Range("E40") "3.53%"
Set-TextValue $ws.Range("G40") "18"

Set-TextValue $ws.Range("D41") "0.006569"
Set-TextValue $ws.Range("E41") "-1.70%"
Set-TextValue $ws.Range("G41") "18"

Set-TextValue $ws.Range("D42") "0.1241"
Set-TextValue $ws.Range("E42") "-11.23%"
Set-TextValue $ws.Range("G42") "18"

Set-TextValue $ws.Range("D43") "0.002003"
Set-TextValue $ws.Range("E43") "-3.88%"
Set-TextValue $ws.Range("G43") "18"

Set-TextValue $ws.Range("D44") "0.01147"
Set-TextValue $ws.Range("E44") "-0.82%"
Set-TextValue $ws.Range("G44") "18"

Set-TextValue $ws.Range("D45") "0.00005476"
Set-TextValue $ws.Range("E45") "-0.81%"
Set-TextValue $ws.Range("G45") "18"

Set-TextValue $ws.Range("E46") "25.93%"
Set-TextValue $ws.Range("G46") "18"

Set-TextValue $ws.Range("D47") "0.01296"
Set-TextValue $ws.Range("E47") "-29.72%"
Set-TextValue $ws.Range("G47") "18"

Set-TextValue $ws.Range("G48") "18"

Set-TextValue $ws.Range("G49") "18"

Set-TextValue $ws.Range("G50") "18"

Set-TextValue $ws.Range("G51") "18"
